$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.612.66"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.622.71"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'214.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'19.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.850.03"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "1.624.17"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "'0.515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "'64.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "26.592.96"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "'230.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.01%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "'4.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").Value = "'9.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "'2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("D25").Value = "'145.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").Value = "'15.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").Value = "1.445.53"
$ws.Range("E33").Value = "  +8.03%  "
$ws.Range("D34").Value = "'3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").Value = "'0.562"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "'0.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").Value = "'5.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").Value = "1.760.69"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'0.929"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.79%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.765"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'62.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").Value = "'88.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").Value = "'1.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "'0.0968"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
